# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Worker "AIDEE REYES BARRIOS" (CC 45373580): the two "Periodo Mora" rows
# (16 & 17) swap order -> E16 becomes 2211, E17 becomes 2212.
#
# Worker "LUZ MARINA BLANCO SEGOVIA" (CC 52401661): the two "Periodo Mora"
# rows (18 & 19) swap order -> E18 becomes 2309, E19 becomes 2310, and their
# "Valor Mora" / "Salario Basico" values are updated accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep these as text (they already render like numbers, e.g. "2211"),
# matching the shared-string ("t=s") type of the original cells.
$ws.Range("E16").Value = "2211"
$ws.Range("E17").Value = "2212"

$ws.Range("E18").Value = "2309"
$ws.Range("F18").Value = 37120
$ws.Range("G18").Value = 1423500

$ws.Range("E19").Value = "2310"
$ws.Range("F19").Value = 46400
$ws.Range("G19").Value = 1423500
